# Moved statistical datasets and results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to make room for the row labels.
$ws.Range("A1").EntireColumn.Insert()

# Insert a new row above row 1 to make room for the header row.
$ws.Range("A1").EntireRow.Insert()

# Header row (row 1) for the former A:D columns, now B:E.
$ws.Range("B1").Value = "Valid"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = "Z"
$ws.Range("E1").Value = "p-value"

# Row labels for the data rows (now rows 2..22, previously rows 1..21).
$labels = @(
  "CyclomaticComplexity(CC) & NbOperators",
  "MaintainabilityIndex & MaintainabilityIndex",
  "NbUniqueOperands & NbUniqueOperands",
  "NbUniqueOperands & EffortToImplement",
  "NbOperands & NbOperands",
  "NbOperands & EffortToImplement",
  "NbUniqueOperators & NbUniqueOperators",
  "NbUniqueOperators & EffortToImplement",
  "NbOperators & NbOperators",
  "ProgramLength & ProgramLength",
  "ProgramLength & EffortToImplement",
  "VocabularySize & VocabularySize",
  "ProgramVolume & ProgramVolume",
  "DifficultyLevel & DifficultyLevel",
  "ProgramLevel & ProgramLevel",
  "EffortToImplement & NbUniqueOperands",
  "EffortToImplement & NbOperands",
  "EffortToImplement & NbUniqueOperators",
  "EffortToImplement & ProgramLength",
  "EffortToImplement & EffortToImplement",
  "TimeToImplement & TimeToImplement"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# Widen the new label column. Columns B:E already carry over the original
# A:D widths untouched by the insert, so they are left alone.
$ws.Columns.Item(1).ColumnWidth = 53.65
